# Repull data, push all data, mean calculation
# Update the dSF column (F) values for several rows to reflect the latest data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -2
$ws.Range("F4").Value = 0
$ws.Range("F6").Value = -4
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -1
$ws.Range("F9").Value = 5
